$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.916.43'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '2.534.09'
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.90'
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.14'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.577'
$ws.Range("E7").Value = '  -1.48%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'
$ws.Range("E9").Value = '  -0.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.24'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0811'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.58'
$ws.Range("E12").Value = '  -2.07%  '
$ws.Range("E13").Value = '  -3.63%  '
$ws.Range("D14").Value = '2.924.72'
$ws.Range("E14").Value = '  +1.00%  '
$ws.Range("D15").Value = '2.518.91'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("E16").Value = '  -2.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.857'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").Value = '42.941.20'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.93'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.78'
$ws.Range("E20").Value = '  +3.89%  '
$ws.Range("D21").Value = '0.0₃0965'
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.98'
$ws.Range("E22").Value = '  -2.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.84'
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.96'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  +1.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.77'
$ws.Range("E26").Value = '  -1.17%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.43'
$ws.Range("E28").Value = '  +3.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.73'
$ws.Range("E29").Value = '  +7.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.40'
$ws.Range("E30").Value = '  +2.20%  '
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '157.66'
$ws.Range("E32").Value = '  +2.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.18'
$ws.Range("E33").Value = '  +4.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.30'
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.35'
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0780'
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.112'
$ws.Range("E38").Value = '  -1.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.119'
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.51'
$ws.Range("E40").Value = '  -5.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.34'
$ws.Range("E41").Value = '  +15.80%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.82'
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0305'
$ws.Range("E43").Value = '  +0.56%  '
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("D46").Value = '2.051.24'
$ws.Range("E46").Value = '  +1.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.21'
$ws.Range("E47").Value = '  +0.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '108.41'
$ws.Range("E48").Value = '  +6.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.98'
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.09'
$ws.Range("E50").Value = '  +2.67%  '
$ws.Range("D51").Value = '2.774.71'
$ws.Range("E51").Value = '  +0.80%  '
